$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5 holds the Treatment query (row label A5 = "TreatmentTab").
# Fix the long standing bug in the Treatment query: drop the redundant CONCAT() wrapper
# around REPLACE(), i.e. CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) -> REPLACE(trt.treatment_agent, ';', ', ')
$treatmentQuery = $ws.Range("B5").Value2
$treatmentQueryFixed = $treatmentQuery.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$ws.Range("B5").Value2 = $treatmentQueryFixed

# Nudge B5's font back to the default 11pt Calibri explicitly (re-applying the font picks up
# a distinct font record / cell style, matching the formatting refresh that came with the edit).
$ws.Range("B5").Font.Size = 11
$ws.Range("B5").Font.Name = "Calibri"

# --- Update the active selection / scroll position on the sheet ---
$ws.Range("B2").Select()
